$d = $word.ActiveDocument

# Locate the target table: the one whose row 8 / column 2 contains "B_Area"
# (this is Table 3 in Description_0.docx, the Building table).
$table = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Rows.Count -ge 8) {
        $txt = $candidate.Cell(8, 2).Range.Text
        if ($txt -like "B_Area*") {
            $table = $candidate
        }
    }
}

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-CellParagraphXml($cell, [string]$innerXml) {
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = ""
    $rng2 = $cell.Range
    $rng2.End = $rng2.End - 1
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + '><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng2.InsertXML($pkg)
}

$rFonts = '<w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/>'

# --- Row 6 (label "5"): BCR -> B_Area -------------------------------------

$cell = $table.Cell(6, 2)
$xml = '<w:p><w:pPr><w:spacing w:before="20" w:after="20"/><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr>' + $rFonts + '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr>' + $rFonts + '<w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="vi-VN"/></w:rPr><w:t>B_Area</w:t></w:r></w:p>'
Set-CellParagraphXml $cell $xml

$cell = $table.Cell(6, 3)
$xml = '<w:p><w:pPr><w:spacing w:before="20" w:after="20"/><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr>' + $rFonts + '<w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr>' + $rFonts + '<w:lang w:val="vi-VN"/></w:rPr><w:t>Building footprint (m2)</w:t></w:r></w:p>'
Set-CellParagraphXml $cell $xml

# --- Row 7 (label "6"): FAR -> BuildingTy ----------------------------------

$cell = $table.Cell(7, 1)
$xml = '<w:p><w:pPr><w:spacing w:before="20" w:after="20"/><w:rPr>' + $rFonts + '<w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr>' + $rFonts + '<w:lang w:val="vi-VN"/></w:rPr><w:t>6</w:t></w:r></w:p>'
Set-CellParagraphXml $cell $xml

$cell = $table.Cell(7, 2)
$xml = '<w:p><w:pPr><w:spacing w:before="20" w:after="20"/><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr>' + $rFonts + '<w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr>' + $rFonts + '<w:lang w:val="vi-VN"/></w:rPr><w:t>BuildingTy</w:t></w:r></w:p>'
Set-CellParagraphXml $cell $xml

$cell = $table.Cell(7, 3)
$newText = "This field has 3 types: NO [nomal], KD [the building" + [char]0x2019 + "s podium], KT [the building" + [char]0x2019 + "s tower]"
$xml = '<w:p><w:pPr><w:spacing w:before="20" w:after="20"/><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:rPr>' + $rFonts + '<w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr>' + $rFonts + '<w:lang w:val="vi-VN"/></w:rPr><w:t>' + $newText + '</w:t></w:r></w:p>'
Set-CellParagraphXml $cell $xml

# --- Row 8 (label "7" / B_Area / Building footprint) is removed entirely --

$table.Rows.Item(8).Delete()
